$p = $ppt.ActivePresentation

# --- 1) Slide 2: trim the run that used to continue with the YOLO sentence ---
$s2 = $p.Slides.Item(2)
$contentShape = $s2.Shapes.Item(2)
$tr2 = $contentShape.TextFrame.TextRange
$run = $tr2.Paragraphs(3).Runs(3)
$run.Text = " use cameras placed at intersections to capture images of each lane. "

# --- 2) Slide 5: add the "Simulator link" textbox under the IOT picture ---
$s5 = $p.Slides.Item(5)

$left   = 3959225 / 12700
$top    = 5725020 / 12700
$width  = 7395315 / 12700
$height = 646331 / 12700

$tb = $s5.Shapes.AddTextbox(1, $left, $top, $width, $height)
$tb.Name = "TextBox 2"
$tb.Fill.Visible = 0
$tb.TextFrame.WordWrap = -1

$prefix = "Simulator link:- "
$url = "https://wokwi.com/projects/373481220669138945"

$tbr = $tb.TextFrame.TextRange
$tbr.Text = $prefix + $url + "`r"

$urlRange = $tbr.Characters($prefix.Length + 1, $url.Length)
$urlRange.ActionSettings(1).Hyperlink.Address = $url

$tb.TextFrame.AutoSize = 1

Write-Output "done"
